# 明日之后共创服数据库 - 副本.xlsx
# Remove the columns that are no longer needed from the "副本" (dungeon)
# table: 副本模式(D), 等级区间(F), 反击难度(G), 建议配置(I), 推荐武器(J).
# Deleting whole columns shifts everything after them to the left, which
# turns the original A:M header row (副本编号, 副本名称, 所在地图, 副本模式,
# 副本类型, 等级区间, 反击难度, 普通或讨伐, 建议配置, 推荐武器, 异变类型)
# into the trimmed A:H header row (副本编号, 副本名称, 所在地图, 副本类型,
# 普通或讨伐, 异变类型).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from right to left so the earlier deletes don't invalidate the
# column indices of the ones still pending.
$ws.Columns.Item(10).Delete()   # J - 推荐武器
$ws.Columns.Item(9).Delete()    # I - 建议配置
$ws.Columns.Item(7).Delete()    # G - 反击难度
$ws.Columns.Item(6).Delete()    # F - 等级区间
$ws.Columns.Item(4).Delete()    # D - 副本模式

# Leave the selection on column F (now "异变类型"), matching the saved
# workbook's UI state after the edit.
$ws.Columns.Item(6).Select() | Out-Null
